$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.929.91'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.508.40'
$ws.Range('E3').Value = '  +0.23%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('E10').Value = '  -1.31%  '
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '2.954.69'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').Value = '58.849.16'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.73'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('E16').Value = '  -1.33%  '
$ws.Range('D17').Value = '2.511.26'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.94'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.63%  '
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('D28').Value = '0.0₃0762'
$ws.Range('E28').Value = '  -2.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.48'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '169.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('E33').Value = '  -5.02%  '
$ws.Range('E34').Value = '  -3.27%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.40'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.16%  '
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('E37').Value = '  -2.96%  '
$ws.Range('E38').Value = '  -2.20%  '
$ws.Range('E39').Value = '  -4.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '281.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.74%  '
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.604'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.50%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.06%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '129.66'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0926'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('E47').Value = '  -2.62%  '
$ws.Range('E48').Value = '  -3.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.24'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.15%  '
$ws.Range('D50').Value = '1.757.96'
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('E51').Value = '  -0.44%  '
